$d = $word.ActiveDocument
$d.Content.Find.Execute("2025-11-08 Saturday", $false, $false, $false, $false, $false, $true, 1, $false, "2025-11-09 Sunday", 2) | Out-Null
$d.Content.Find.Execute("848×7=5936", $false, $false, $false, $false, $false, $true, 1, $false, "887×3=2661", 2) | Out-Null
$d.Content.Find.Execute("857×5=4285", $false, $false, $false, $false, $false, $true, 1, $false, "670×6=4020", 2) | Out-Null
$d.Content.Find.Execute("684×5=3420", $false, $false, $false, $false, $false, $true, 1, $false, "636×8=5088", 2) | Out-Null
$d.Content.Find.Execute("665×7=4655", $false, $false, $false, $false, $false, $true, 1, $false, "979×8=7832", 2) | Out-Null
$d.Content.Find.Execute("176×5=880", $false, $false, $false, $false, $false, $true, 1, $false, "124×9=1116", 2) | Out-Null
$d.Content.Find.Execute("796×2=1592", $false, $false, $false, $false, $false, $true, 1, $false, "122×3=366", 2) | Out-Null
$d.Content.Find.Execute("592×7=4144", $false, $false, $false, $false, $false, $true, 1, $false, "525×3=1575", 2) | Out-Null
$d.Content.Find.Execute("661×9=5949", $false, $false, $false, $false, $false, $true, 1, $false, "175×3=525", 2) | Out-Null
$d.Content.Find.Execute("972×9=8748", $false, $false, $false, $false, $false, $true, 1, $false, "587×6=3522", 2) | Out-Null
$d.Content.Find.Execute("657×4=2628", $false, $false, $false, $false, $false, $true, 1, $false, "519×6=3114", 2) | Out-Null
$d.Content.Find.Execute("191×8=1528", $false, $false, $false, $false, $false, $true, 1, $false, "618×5=3090", 2) | Out-Null
$d.Content.Find.Execute("524×2=1048", $false, $false, $false, $false, $false, $true, 1, $false, "219×6=1314", 2) | Out-Null
$d.Content.Find.Execute("405×2=810", $false, $false, $false, $false, $false, $true, 1, $false, "934×3=2802", 2) | Out-Null
$d.Content.Find.Execute("740×2=1480", $false, $false, $false, $false, $false, $true, 1, $false, "325×2=650", 2) | Out-Null
$d.Content.Find.Execute("611×9=5499", $false, $false, $false, $false, $false, $true, 1, $false, "251×6=1506", 2) | Out-Null
$d.Content.Find.Execute("652×3=1956", $false, $false, $false, $false, $false, $true, 1, $false, "225×3=675", 2) | Out-Null
$d.Content.Find.Execute("394×9=3546", $false, $false, $false, $false, $false, $true, 1, $false, "765×6=4590", 2) | Out-Null
$d.Content.Find.Execute("853×3=2559", $false, $false, $false, $false, $false, $true, 1, $false, "319×7=2233", 2) | Out-Null
$d.Content.Find.Execute("641×2=1282", $false, $false, $false, $false, $false, $true, 1, $false, "148×4=592", 2) | Out-Null
$d.Content.Find.Execute("475×5=2375", $false, $false, $false, $false, $false, $true, 1, $false, "258×8=2064", 2) | Out-Null
$d.Content.Find.Execute("458×9=4122", $false, $false, $false, $false, $false, $true, 1, $false, "706×6=4236", 2) | Out-Null
$d.Content.Find.Execute("717×6=4302", $false, $false, $false, $false, $false, $true, 1, $false, "747×4=2988", 2) | Out-Null
$d.Content.Find.Execute("372×2=744", $false, $false, $false, $false, $false, $true, 1, $false, "443×5=2215", 2) | Out-Null
$d.Content.Find.Execute("479×6=2874", $false, $false, $false, $false, $false, $true, 1, $false, "921×9=8289", 2) | Out-Null
$d.Content.Find.Execute("152×7=1064", $false, $false, $false, $false, $false, $true, 1, $false, "962×2=1924", 2) | Out-Null
